# Add "status_label" column: a string version of the "statut" emoji column.
# Insert a new column before column B (NCTId), shifting old B..I to C..J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Insert()

# New header for the inserted column B: same look as the other header cells
# (bold font, thin border, centered/top-aligned) used by A1/C1/etc.
$ws.Cells.Item(1, 2).Value = "status_label"
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160
$b1.Borders.LineStyle = 1

# Map the emoji in column A to its French colour-name label in new column B.
$statusMap = @{
    "🟥" = "rouge"
    "🟧" = "orange"
    "🟨" = "jaune"
    "🟩" = "vert"
}

for ($r = 2; $r -le 14; $r++) {
    $statut = $ws.Cells.Item($r, 1).Value2
    $label = $statusMap[$statut]
    $ws.Cells.Item($r, 2).Value = $label
}

Write-Host "status_label column added"
